# Add missing supply prefixes to the SupplyLookupMappings sheet.
# Two new prefix -> supply rows are appended at the bottom of the table:
#   VL-MG -> Valero
#   GMK   -> Growmark

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SupplyLookupMappings")
$ws.Activate()

$ws.Range("A33").Value = "VL-MG"
$ws.Range("B33").Value = "Valero"

$ws.Range("A34").Value = "GMK"
$ws.Range("B34").Value = "Growmark"

# Match the author's final cursor position after the edit.
[void]$ws.Range("D5").Select()
